# Updates cryptos list values (Price / Volume(1h) columns, plus a row-20/21
# coin swap) to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a numeric-looking text cell (Price / Volume(1h))
# without losing its text formatting (e.g. '62.922.46', '1.00', '  -0.25%  ').
# Force the cell to Text before assigning, then restore the default style so
# no stray number-format/style diff is introduced.
function Set-TextCell([string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $val
    $cell.Style = 'Normal'
}

Set-TextCell 'D2' '62.922.46'
Set-TextCell 'E2' '  -0.25%  '
Set-TextCell 'D3' '2.543.00'
Set-TextCell 'E3' '  +3.57%  '
Set-TextCell 'E4' '  +0.10%  '
Set-TextCell 'D5' '568.97'
Set-TextCell 'E5' '  +0.83%  '
Set-TextCell 'D6' '145.63'
Set-TextCell 'E6' '  +2.34%  '
Set-TextCell 'E7' '  +0.07%  '
Set-TextCell 'E8' '  +0.21%  '
Set-TextCell 'D9' '2.540.64'
Set-TextCell 'E9' '  +3.51%  '
Set-TextCell 'E10' '  +0.29%  '
Set-TextCell 'D11' '5.51'
Set-TextCell 'E11' '  -2.64%  '
Set-TextCell 'E12' '  -0.27%  '
Set-TextCell 'E13' '  -0.54%  '
Set-TextCell 'D14' '27.23'
Set-TextCell 'E14' '  +0.58%  '
Set-TextCell 'D15' '2.997.38'
Set-TextCell 'E15' '  +3.64%  '
Set-TextCell 'D16' '62.865.53'
Set-TextCell 'E16' '  +0.07%  '
Set-TextCell 'D17' '0.0000143'
Set-TextCell 'E17' '  +1.35%  '
Set-TextCell 'D18' '2.548.53'
Set-TextCell 'E18' '  +3.93%  '
Set-TextCell 'D19' '11.30'
Set-TextCell 'E19' '  +0.78%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 'D20' '334.07'
Set-TextCell 'E20' '  -1.56%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D21' '4.33'
Set-TextCell 'E21' '  +1.38%  '
Set-TextCell 'D22' '6.78'
Set-TextCell 'E22' '  +0.52%  '
Set-TextCell 'D23' '1.00'
Set-TextCell 'E23' '  +0.13%  '
Set-TextCell 'E24' '  -0.46%  '
Set-TextCell 'E25' '  -0.78%  '
Set-TextCell 'E26' '  +6.92%  '
Set-TextCell 'E27' '  +0.09%  '
Set-TextCell 'D28' '1.48'
Set-TextCell 'E28' '  +3.49%  '
Set-TextCell 'D29' '8.35'
Set-TextCell 'E29' '  +3.62%  '
Set-TextCell 'D30' '7.29'
Set-TextCell 'E30' '  +8.32%  '
Set-TextCell 'E31' '  +3.30%  '
Set-TextCell 'E32' '  +0.68%  '
Set-TextCell 'D33' '175.26'
Set-TextCell 'E33' '  -0.68%  '
Set-TextCell 'E34' '  +2.55%  '
Set-TextCell 'D35' '408.39'
Set-TextCell 'E35' '  +5.86%  '
Set-TextCell 'D36' '0.398'
Set-TextCell 'E36' '  +0.47%  '
Set-TextCell 'D37' '18.96'
Set-TextCell 'E37' '  +1.22%  '
Set-TextCell 'E38' '  +0.01%  '
Set-TextCell 'D39' '4.34'
Set-TextCell 'E39' '  +0.44%  '
Set-TextCell 'E40' '  +0.80%  '
Set-TextCell 'E41' '  +0.08%  '
Set-TextCell 'D42' '39.62'
Set-TextCell 'E42' '  -0.89%  '
Set-TextCell 'D43' '151.85'
Set-TextCell 'E43' '  +1.66%  '
Set-TextCell 'E44' '  +1.59%  '
Set-TextCell 'D45' '20.74'
Set-TextCell 'E45' '  +1.27%  '
Set-TextCell 'E46' '  +0.98%  '
Set-TextCell 'D47' '0.0530'
Set-TextCell 'E47' '  +3.12%  '
Set-TextCell 'D48' '0.0962'
Set-TextCell 'E48' '  +0.16%  '
Set-TextCell 'E49' '  +4.05%  '
Set-TextCell 'D50' '18.25'
Set-TextCell 'E50' '  +2.10%  '
Set-TextCell 'D51' '1.73'
Set-TextCell 'E51' '  -1.98%  '
